$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data row down to make room for a new header row.
$ws.Rows.Item(1).Insert()

# ---- Header row ----
$ws.Range("A1").NumberFormat = "@"
$ws.Range("D1").NumberFormat = "@"
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Value"
$ws.Range("C1").Value = "weight"
$ws.Range("D1").Value = "baseid"
$ws.Range("E1").Value = "item hp"
$ws.Range("F1").Value = "damage rating"

# ---- Row 2 (pre-existing item) keeps its values; just apply the Text format ----
$ws.Range("A2").NumberFormat = "@"
$ws.Range("D2").NumberFormat = "@"

# ---- New armor rows ----
$ws.Range("A3").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("A3").Value = "Chinese Jumpsuit"
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "00078646"
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = 1

$ws.Range("A4").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("A4").Value = "Armored Vault 101 Jumpsuit"
$ws.Range("B4").Value = 180
$ws.Range("C4").Value = 15
$ws.Range("D4").Value = "00034121"
$ws.Range("E4").Value = 100
$ws.Range("F4").Value = 12

$ws.Range("A5").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("A5").Value = "Child's Vault 101 Jumpsuit"
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = "000340f2"
$ws.Range("E5").Value = 100
$ws.Range("F5").Value = 1

$ws.Range("A6").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("A6").Value = "Dad's Wasteland Outfit"
$ws.Range("B6").Value = 6
$ws.Range("C6").Value = 20
$ws.Range("D6").Value = "00079f09"
$ws.Range("E6").Value = 100
$ws.Range("F6").Value = 2

$ws.Range("A7").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("A7").Value = "Modified Utility Jumpsuit"
$ws.Range("B7").Value = 30
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = "0007c17c"
$ws.Range("E7").Value = 100
$ws.Range("F7").Value = 1

$ws.Range("A8").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("A8").Value = "Tunnel Snake Outfit"
$ws.Range("B8").Value = 8
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = "0002042e"
$ws.Range("E8").Value = 100
$ws.Range("F8").Value = 4

# ---- Cosmetics: widen the name/value columns (content grew) and set print orientation ----
$ws.Columns.Item(1).ColumnWidth = 25.592447916666668
$ws.Columns.Item(2).ColumnWidth = 5.307291666666667
$ws.PageSetup.Orientation = 1

# ---- Selection / cursor position ----
$null = $ws.Range("F16").Select()
